$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 237, shifting existing rows 237:267 down to 238:268
$ws.Rows("237").Insert()

# Populate the newly inserted row 237 with the new record's data
$ws.Range("A237").Value = 8
$ws.Range("B237").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C237").Value = 'Coquimbo'
$ws.Range("D237").Value = 45124
$ws.Range("D237").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E237").Value = 4
$ws.Range("F237").Value = 100112001
$ws.Range("G237").Value = 'Berenjena'
$ws.Range("H237").Value = 'Sin especificar'
$ws.Range("I237").Value = 'Primera'
$ws.Range("J237").Value = 400
$ws.Range("K237").Value = 8500
$ws.Range("L237").Value = 9000
$ws.Range("M237").Value = 8750
$ws.Range("N237").Value = '$/caja 50 unidades'
$ws.Range("O237").Value = 'Región de Arica y Parinacota'
$ws.Range("P237").Value = 175
$ws.Range("Q237").Value = 50
$ws.Range("R237").Value = 'Hortaliza'
